# Updates the cryptos list with refreshed price / volume(1h) data.
# Generated to match the commit "Updated cryptos list on Fri Jul 14 10:31:33 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.146.46'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '1.990.53'
$ws.Range("E3").Value = '  +5.76%  '
$ws.Range("D4").Value = '''0.9993'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''0.7915'
$ws.Range("E5").Value = '  +67.46%  '
$ws.Range("D6").Value = '''254.00'
$ws.Range("E6").Value = '  +3.06%  '
$ws.Range("D7").Value = '''0.9991'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '''0.3497'
$ws.Range("E8").Value = '  +21.24%  '
$ws.Range("D9").Value = '''28.07'
$ws.Range("E9").Value = '  +26.99%  '
$ws.Range("D10").Value = '''0.06993'
$ws.Range("E10").Value = '  +7.05%  '
$ws.Range("D11").Value = '''0.8451'
$ws.Range("E11").Value = '  +9.87%  '
$ws.Range("D12").Value = '''0.08188'
$ws.Range("E12").Value = '  +4.65%  '
$ws.Range("D13").Value = '1.992.38'
$ws.Range("E13").Value = '  +6.01%  '
$ws.Range("D14").Value = '''100.28'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("E15").Value = '  +6.74%  '
$ws.Range("D16").Value = '''15.36'
$ws.Range("E16").Value = '  +16.53%  '
$ws.Range("D17").Value = '''272.97'
$ws.Range("E17").Value = '  -4.22%  '
$ws.Range("D18").Value = '31.152.29'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").Value = '''5.860'
$ws.Range("E19").Value = '  +9.37%  '
$ws.Range("D20").Value = '''0.000007894'
$ws.Range("E20").Value = '  +4.99%  '
$ws.Range("D21").Value = '2.259.50'
$ws.Range("E21").Value = '  +6.34%  '
$ws.Range("D22").Value = '''0.9995'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '''0.9991'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '''7.062'
$ws.Range("E24").Value = '  +10.42%  '
$ws.Range("D25").Value = '''10.06'
$ws.Range("E25").Value = '  +10.16%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.1505'
$ws.Range("E26").Value = '  +55.23%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''164.61'
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D28").Value = '''19.86'
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("D29").Value = '''2.321'
$ws.Range("E29").Value = '  +21.30%  '
$ws.Range("D30").Value = '''1.596'
$ws.Range("E30").Value = '  +6.09%  '
$ws.Range("D31").Value = '''1.360'
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("D32").Value = '''4.585'
$ws.Range("E32").Value = '  +7.65%  '
$ws.Range("D33").Value = '''4.402'
$ws.Range("E33").Value = '  +4.94%  '
$ws.Range("D34").Value = '''0.05223'
$ws.Range("E34").Value = '  +8.00%  '
$ws.Range("D35").Value = '''1.229'
$ws.Range("E35").Value = '  +8.84%  '
$ws.Range("D36").Value = '''0.7770'
$ws.Range("E36").Value = '  +11.68%  '
$ws.Range("D37").Value = '''2.758'
$ws.Range("D38").Value = '''0.02003'
$ws.Range("E38").Value = '  +4.69%  '
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").Value = '''6.629'
$ws.Range("E40").Value = '  +5.62%  '
$ws.Range("D41").Value = '''79.31'
$ws.Range("E41").Value = '  +3.81%  '
$ws.Range("D42").Value = '''0.4664'
$ws.Range("E42").Value = '  +9.59%  '
$ws.Range("D43").Value = '''2.123'
$ws.Range("E43").Value = '  +7.34%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''104.80'
$ws.Range("E44").Value = '  +3.14%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''0.8472'
$ws.Range("E45").Value = '  +2.16%  '
$ws.Range("D46").Value = '''0.9992'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = '''7.674'
$ws.Range("E47").Value = '  +9.12%  '
$ws.Range("D48").Value = '''9.850'
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''36.82'
$ws.Range("E49").Value = '  +4.81%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '''0.4296'
$ws.Range("E50").Value = '  +8.80%  '
$ws.Range("D51").Value = '''1.527'
$ws.Range("E51").Value = '  +13.08%  '
